# Business Summary slide (slide 2): tighten the sub-headline copy/size and
# trim the first bullet, matching the updated bodyPr/size on the bullet box.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Shape "Google Shape;66;p14" (body placeholder, idx=1) ---------------
$headline = $s.Shapes.Item(2)
$headlineRange = $headline.TextFrame.TextRange

# Left-align the paragraph (was centered).
$headlineRange.Paragraphs(1, 1).ParagraphFormat.Alignment = 1

# Update the run text and shrink the font from 24pt to 20pt.
$headlineRun = $headlineRange.Runs(1, 1)
$headlineRun.Text = "Analysis of industry data reveal opportunities for success through investment in: "
$headlineRun.Font.Size = 20

# --- Shape "TextBox 4" (bullet list) --------------------------------------
$bullets = $s.Shapes.Item(3)

# Narrow the box to match the placeholders above (8744702 EMU -> 8520600 EMU)
$bullets.Width = 8520600 / 12700

# Let the text wrap inside the (now narrower) box instead of forcing one line.
$bullets.TextFrame.WordWrap = -1

# Trim the first bullet's text.
$bulletsRange = $bullets.TextFrame.TextRange
$firstBulletRun = $bulletsRange.Runs(1, 1)
$firstBulletRun.Text = "Popular and highly rated genres"

# Re-pin the auto-fit height (unchanged by the source edit) after the text
# and wrap/width changes above, which would otherwise relayout the box.
$bullets.Height = (2236510 / 12700) + 0.00001
